$wb = $excel.ActiveWorkbook

# Rename worksheets
$wb.Worksheets.Item(1).Name = "iCC389"
$wb.Worksheets.Item(3).Name = "iCC470"
$wb.Worksheets.Item(4).Name = "iCC651"

# Sheet1 value updates
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 0.1334525431033355
$ws.Range("B7").Value = 0.149997412972611
$ws.Range("B8").Value = 0.1169076732340626
$ws.Range("B10").Value = 0.03418332388769043
$ws.Range("B11").Value = 0.1499974129726072
$ws.Range("B12").Value = 0.1499974129726087

# Sheet2 value updates
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 0.9599355516605008
$ws.Range("B3").Value = 0.9599355516604995
$ws.Range("B7").Value = 0.8009057833985511
$ws.Range("B8").Value = 0.800905783398551
$ws.Range("B10").Value = 0.9599355516604984

# Sheet3 value updates
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 0.7515130259435367
$ws.Range("B3").Value = 0.8028664160496785
$ws.Range("B4").Value = 0.8028664160496776
$ws.Range("B5").Value = 0.8028664160496781
$ws.Range("B6").Value = 0.766326397683184
$ws.Range("B8").Value = 0.7001596358373952
$ws.Range("B10").Value = 0.4433926853066865
$ws.Range("B12").Value = 0.4536865487071338
$ws.Range("B13").Value = 0.8028664160496781
$ws.Range("B14").Value = 0.8028664160496787
$ws.Range("B15").Value = 0.8028664160496783
$ws.Range("B16").Value = 0.5974528556251114
$ws.Range("B17").Value = 0.8078765028893019
$ws.Range("B18").Value = 0.3741968213335981

# Sheet4 value updates
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 0.332754769563311
$ws.Range("B3").Value = 0.05043421551599262
$ws.Range("B4").Value = 0.3167956982730797
$ws.Range("B5").Value = 0.3167956982730806
$ws.Range("B6").Value = 0.3092868722039628
$ws.Range("B7").Value = 0.3167956982730808
$ws.Range("B8").Value = 0.3092868722039628
$ws.Range("B9").Value = 0.3092868722039628
$ws.Range("B10").Value = 0.2050821992414801
$ws.Range("B11").Value = 0.3327547695633093
$ws.Range("B12").Value = 0.2050821992414801
$ws.Range("B13").Value = 0.3327547695633112
$ws.Range("B14").Value = 0.3327547695633093
$ws.Range("B15").Value = 0.2689184844023939
$ws.Range("B16").Value = 0.2884930403580109
$ws.Range("B17").Value = 0.05043421551599279
$ws.Range("B18").Value = 0.3092868722039629
$ws.Range("B19").Value = 0.3167956982730797
$ws.Range("B20").Value = 0.05043421551599279
$ws.Range("B21").Value = 0.1534788460993879
